$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: insert new "lang_code" column before "code" (column A),
# shifting the rest of the former headers one column to the right.
# ---------------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "lang_code"
$ws.Cells.Item(1,2).Value = "code"
$ws.Cells.Item(1,3).Value = "name"
$ws.Cells.Item(1,4).Value = "day_seq"
$ws.Cells.Item(1,5).Value = "is_global_working"
$ws.Cells.Item(1,6).Value = "is_active"

# ---------------------------------------------------------------------
# Data rows 2-15: lang_code, code, name, day_seq, is_global_working, is_active
# ---------------------------------------------------------------------
$rows = @(
  @("eng", 101, "SUN", 1, $false, $true),
  @("eng", 102, "MON", 2, $true,  $true),
  @("eng", 103, "TUE", 3, $true,  $true),
  @("eng", 104, "WED", 4, $true,  $true),
  @("eng", 105, "THU", 5, $true,  $true),
  @("eng", 106, "FRI", 6, $true,  $true),
  @("eng", 107, "SAT", 7, $false, $true),
  @("fra", 101, "DIM", 1, $false, $true),
  @("fra", 102, "LUN", 2, $true,  $true),
  @("fra", 103, "MAR", 3, $true,  $true),
  @("fra", 104, "MER", 4, $true,  $true),
  @("fra", 105, "JEU", 5, $true,  $true),
  @("fra", 106, "VEN", 6, $true,  $true),
  @("fra", 107, "SAM", 7, $false, $true)
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r,1).Value = $row[0]
  $ws.Cells.Item($r,2).Value = $row[1]
  $ws.Cells.Item($r,3).Value = $row[2]
  $ws.Cells.Item($r,4).Value = $row[3]
  $ws.Cells.Item($r,5).Value = $row[4]
  $ws.Cells.Item($r,6).Value = $row[5]
  $r = $r + 1
}

# ---------------------------------------------------------------------
# Apply the header's (bold / bordered / centered) style to every cell
# in the new "lang_code" column, matching the style already used on A1.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A2:A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
